$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.039.26"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.554.49"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'288.90"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.3939"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("D8").Value = "'0.3231"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'42.79"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").Value = "'0.07331"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'1.099"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "'19.01"
$ws.Range("E13").Value = "  -6.74%  "
$ws.Range("D14").Value = "'5.646"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "'0.00001146"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "'6.697"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "1.558.99"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "'0.06594"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'83.96"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "'0.9997"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'6.347"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'15.81"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'11.30"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "22.127.84"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'2.342"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").Value = "'2.463"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("D27").Value = "'148.65"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'18.78"
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").Value = "'4.869"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "1.732.17"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'119.72"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("D32").Value = "'1.061"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'5.726"
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("D34").Value = "'0.08366"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "'9.266"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").Value = "'1.614"
$ws.Range("E36").Value = "  -13.16%  "
$ws.Range("D37").Value = "'0.06233"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'0.02274"
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("D39").Value = "'5.167"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "'1.221"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "'0.2072"
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("D42").Value = "'0.9981"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'10.70"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").Value = "'0.5848"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "'13.23"
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").Value = "'3.742"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'0.5613"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").Value = "'1.909"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("D49").Value = "'117.96"
$ws.Range("E49").Value = "  -4.28%  "
$ws.Range("D50").Value = "'1.142"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").Value = "'0.06854"
$ws.Range("E51").Value = "  -3.18%  "
